$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; existing rows 4-36 shift down to 5-37.
$ws.Rows.Item(4).EntireRow.Insert()

# Populate the newly inserted row 4 with the new weekly data point.
$ws.Range("A4").Value = 9
$ws.Range("B4").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44503
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 100112029
$ws.Range("G4").Value = "Orégano"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 16
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("N4").Value = "`$/docena de atados"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 2833
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = "Hortaliza"
